$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("EMN")

# Row 4 - Inventory
$ws.Range("C4").Value = 1379000000.0
$ws.Range("D4").Value = 1338000000.0
$ws.Range("E4").Value = 1419000000.0
$ws.Range("F4").Value = 1659000000.0
$ws.Range("G4").Value = 1662000000.0

# Row 14 - Accounts Payable
$ws.Range("C14").Value = 799000000.0
$ws.Range("D14").Value = 663000000.0
$ws.Range("E14").Value = 525000000.0
$ws.Range("F14").Value = 770000000.0
$ws.Range("G14").Value = 890000000.0

# Row 20 - Long Term Tax Liability (Deferred)
$ws.Range("C20").Value = 848000000.0
$ws.Range("D20").Value = 928000000.0
$ws.Range("E20").Value = 932000000.0
$ws.Range("F20").Value = 924000000.0
$ws.Range("G20").Value = 849000000.0
